$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "FBA79697"
$ws.Range("B26").Value = "TD310"
$ws.Range("C26").Value = "B0CQX4K9P5"
$ws.Range("D26").Value = "B0CQX4K9P5"
$ws.Range("P26").Value = 53
$ws.Range("T26").Value = 68720.24

$ws.Range("A27").Value = "FBA77113"
$ws.Range("B27").Value = "K1"
$ws.Range("C27").Value = "B01ISNU3X4"
$ws.Range("D27").Value = "B01ISNU3X4"
$ws.Range("P27").Value = 44
$ws.Range("T27").Value = 59911.66

$ws.Range("A28").Value = "FBA79113"
$ws.Range("B28").Value = "TC310"
$ws.Range("C28").Value = "B0BTCXQQ6M"
$ws.Range("D28").Value = "B0BTCXQQ6M"
$ws.Range("P28").Value = 21
$ws.Range("T28").Value = 37524.55

$ws.Range("A29").Value = "FBA79260"
$ws.Range("B29").Value = "G11"
$ws.Range("C29").Value = "B07GVGMW59"
$ws.Range("D29").Value = "B07GVGMW59"
$ws.Range("P29").Value = 7
$ws.Range("T29").Value = 14522.06

$ws.Range("A30").Value = "FBA79696"
$ws.Range("B30").Value = "TD310+"
$ws.Range("C30").Value = "B0CQX3VB1R"
$ws.Range("D30").Value = "B0CQX3VB1R"
$ws.Range("P30").Value = 10
$ws.Range("T30").Value = 14483.07

$ws.Range("A31").Value = "FBA79114"
$ws.Range("B31").Value = "TC310+"
$ws.Range("C31").Value = "B0CCV74CL7"
$ws.Range("D31").Value = "B0CCV74CL7"
$ws.Range("P31").Value = 4
$ws.Range("T31").Value = 9691.52

$ws.Range("A32").Value = "FBA77111"
$ws.Range("B32").Value = "TC30"
$ws.Range("C32").Value = "B08CVP2HXP"
$ws.Range("D32").Value = "B08CVP2HXP"
$ws.Range("P32").Value = 5
$ws.Range("T32").Value = 9334.74

$ws.Range("A33").Value = "FBA77117"
$ws.Range("B33").Value = "S20"
$ws.Range("C33").Value = "B078WNW4YW"
$ws.Range("D33").Value = "B078WNW4YW"
$ws.Range("P33").Value = 4
$ws.Range("T33").Value = 8594.92

$ws.Range("A34").Value = "FBA79116"
$ws.Range("B34").Value = "TC-777 PRO"
$ws.Range("C34").Value = "B0BYHHSLPC"
$ws.Range("D34").Value = "B0BYHHSLPC"
$ws.Range("P34").Value = 5
$ws.Range("T34").Value = 8131.35

$ws.Range("A35").Value = "FBA77106"
$ws.Range("B35").Value = "T20"
$ws.Range("C35").Value = "B082W4B7SX"
$ws.Range("D35").Value = "B082W4B7SX"
$ws.Range("P35").Value = 4
$ws.Range("T35").Value = 7816.97

$ws.Range("A36").Value = "FBA77101"
$ws.Range("B36").Value = "TC-777"
$ws.Range("C36").Value = "B07WLWN2ZT"
$ws.Range("D36").Value = "B07WLWN2ZT"
$ws.Range("P36").Value = 3
$ws.Range("T36").Value = 5311.02

$ws.Range("A37").Value = "FBA77105"
$ws.Range("B37").Value = "T30"
$ws.Range("C37").Value = "B089FVQD3Z"
$ws.Range("D37").Value = "B089FVQD3Z"
$ws.Range("P37").Value = 2
$ws.Range("T37").Value = 4786.44

$ws.Range("A38").Value = "FBA77114"
$ws.Range("B38").Value = "TC-2030"
$ws.Range("C38").Value = "B07TSN2H9D"
$ws.Range("D38").Value = "B07TSN2H9D"
$ws.Range("P38").Value = 1
$ws.Range("T38").Value = 3643.22

$ws.Range("A39").Value = "FBA77110"
$ws.Range("B39").Value = "TM20"
$ws.Range("C39").Value = "B08NDB5NWP"
$ws.Range("D39").Value = "B08NDB5NWP"
$ws.Range("P39").Value = 1
$ws.Range("T39").Value = 2422.88

$ws.Range("A40").Value = "FBA79574"
$ws.Range("B40").Value = "TC30S"
$ws.Range("C40").Value = "B0B4WTHLX5"
$ws.Range("D40").Value = "B0B4WTHLX5"
$ws.Range("P40").Value = 0
$ws.Range("T40").Value = 0

$ws.Range("A41").Value = "FBA79112"
$ws.Range("B41").Value = "T90"
$ws.Range("C41").Value = "B0BRKDLXCD"
$ws.Range("D41").Value = "B0BRKDLXCD"
$ws.Range("P41").Value = 0
$ws.Range("T41").Value = 0

$ws.Range("A42").Value = "FBA79111"
$ws.Range("B42").Value = "TD510"
$ws.Range("C42").Value = "B0BRKFP94K"
$ws.Range("D42").Value = "B0BRKFP94K"
$ws.Range("P42").Value = 0
$ws.Range("T42").Value = 0

$ws.Range("T1").Select()
